# Update the 시간대_3 sheet: column A labels are reassigned (shared-string
# reorder in the source diff resulted in different labels lining up with
# the existing rows) and column B counts are refreshed with new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => (label text for column A, new count for column B)
$rows = @{
    2  = @("11", 88)
    3  = @("24", 85)
    4  = @("13", 67)
    5  = @("09", 65)
    6  = @("08", 62)
    7  = @("10", 60)
    8  = @("17", 55)
    9  = @("16", 54)
    10 = @("14", 51)
    11 = @("18", 42)
    12 = @("15", 38)
    13 = @("07", 36)
    14 = @("19", 32)
    15 = @("20", 30)
    16 = @("22", 30)
    17 = @("23", 29)
    18 = @("05", 27)
    19 = @("06", 26)
    20 = @("12", 25)
    21 = @("02", 20)
    22 = @("21", 19)
    23 = @("01", 18)
    24 = @("04", 14)
    25 = @("03", 9)
}

# Column A holds the hour label as text (values like "08", "09" must keep
# their leading zero, so force the whole label range to text before writing
# the new values - otherwise Excel would auto-convert "08" to the number 8).
$ws.Range("A2:A25").NumberFormat = "@"

foreach ($r in $rows.Keys) {
    $label = $rows[$r][0]
    $count = $rows[$r][1]

    $ws.Cells.Item($r, 1).Value = $label

    # Column B holds the numeric count.
    $ws.Cells.Item($r, 2).Value = $count
}
